$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three more fields were appended to the export: old/new identifier codes
# and a status flag. Add their headers in the next free columns (D:F),
# which extends the sheet's used range from A1:C12 to A1:F12.
$ws.Range("D1").Value = "ORG_FAC_IDENOLD"
$ws.Range("E1").Value = "ORG_FAC_IDENNEW"
$ws.Range("F1").Value = "ORG_FAC_STATUS"

# Reflect the saved selection state of the edited file (cursor left on
# the new status column, row 5).
$ws.Range("F5").Select()
